# maj nom tbl JE
# Rename "tblVersion" / "tblJeuSemblable" sheets to follow the "B " (base
# table) naming convention used by the other child tables of tblJeu.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("tblVersion").Name = "B tblVersion"
$wb.Worksheets.Item("tblJeuSemblable").Name = "B tblJeuSemblable"

# Add a new entry on the "B tblPlateformeJeu" sheet and leave it as the
# active/selected sheet+cell, matching the refreshed workbook view.
$wsPlateformeJeu = $wb.Worksheets.Item("B tblPlateformeJeu")
$wsPlateformeJeu.Range("F54").Value = "s"
$wsPlateformeJeu.Activate()
$wsPlateformeJeu.Range("F54").Select()
